$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '64.305.63'
$ws.Range('E2').NumberFormat = '@'
$ws.Cells.Item(2, 5).Value = '  +0.05%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '3.503.98'
$ws.Range('E3').NumberFormat = '@'
$ws.Cells.Item(3, 5).Value = '  -0.29%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Cells.Item(4, 5).Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '587.14'
$ws.Range('E5').NumberFormat = '@'
$ws.Cells.Item(5, 5).Value = '  +0.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '134.62'
$ws.Range('E6').NumberFormat = '@'
$ws.Cells.Item(6, 5).Value = '  +1.04%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Cells.Item(7, 5).Value = '  -0.04%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Cells.Item(8, 5).Value = '  -0.65%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Cells.Item(9, 5).Value = '  +1.03%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Cells.Item(10, 5).Value = '  +0.66%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Cells.Item(11, 5).Value = '  +1.63%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '4.103.71'
$ws.Range('E12').NumberFormat = '@'
$ws.Cells.Item(12, 5).Value = '  +0.08%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Cells.Item(13, 5).Value = '  +1.11%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Cells.Item(14, 5).Value = '  +0.78%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '3.502.30'
$ws.Range('E15').NumberFormat = '@'
$ws.Cells.Item(15, 5).Value = '  -0.16%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '26.09'
$ws.Range('E16').NumberFormat = '@'
$ws.Cells.Item(16, 5).Value = '  -6.02%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '64.321.86'
$ws.Range('E17').NumberFormat = '@'
$ws.Cells.Item(17, 5).Value = '  +0.55%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Cells.Item(18, 5).Value = '  -2.29%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Cells.Item(19, 5).Value = '  +1.21%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Cells.Item(20, 5).Value = '  -4.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '393.11'
$ws.Range('E21').NumberFormat = '@'
$ws.Cells.Item(21, 5).Value = '  +1.79%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Cells.Item(22, 5).Value = '  -1.46%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '3.645.11'
$ws.Range('D24').NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '74.30'
$ws.Range('E24').NumberFormat = '@'
$ws.Cells.Item(24, 5).Value = '  +1.29%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Cells.Item(25, 5).Value = '  -0.08%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Cells.Item(26, 5).Value = '  +2.09%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Cells.Item(27, 5).Value = '  -0.93%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '7.56'
$ws.Range('E28').NumberFormat = '@'
$ws.Cells.Item(28, 5).Value = '  -0.52%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Cells.Item(29, 5).Value = '  +0.10%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Cells.Item(30, 5).Value = '  -0.52%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Cells.Item(31, 5).Value = '  -1.38%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '1.48'
$ws.Range('E32').NumberFormat = '@'
$ws.Cells.Item(32, 5).Value = '  -6.43%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '3.526.89'
$ws.Range('E33').NumberFormat = '@'
$ws.Cells.Item(33, 5).Value = '  +0.23%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Cells.Item(34, 5).Value = '  +0.02%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Cells.Item(35, 5).Value = '  +2.56%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '23.55'
$ws.Range('E36').NumberFormat = '@'
$ws.Cells.Item(36, 5).Value = '  -1.24%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Cells.Item(37, 5).Value = '  -2.58%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '1.58'
$ws.Range('E38').NumberFormat = '@'
$ws.Cells.Item(38, 5).Value = '  -1.11%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '6.92'
$ws.Range('E39').NumberFormat = '@'
$ws.Cells.Item(39, 5).Value = '  -0.98%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '162.76'
$ws.Range('E40').NumberFormat = '@'
$ws.Cells.Item(40, 5).Value = '  +0.75%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Cells.Item(41, 5).Value = '  -2.60%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Cells.Item(42, 5).Value = '  -1.24%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '25.85'
$ws.Range('E43').NumberFormat = '@'
$ws.Cells.Item(43, 5).Value = '  -1.09%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Cells.Item(44, 5).Value = '  +0.21%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '41.87'
$ws.Range('E45').NumberFormat = '@'
$ws.Cells.Item(45, 5).Value = '  +0.44%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Cells.Item(46, 5).Value = '  -0.12%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Cells.Item(47, 5).Value = '  -4.47%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Cells.Item(48, 5).Value = '  -0.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '2.465.92'
$ws.Range('E49').NumberFormat = '@'
$ws.Cells.Item(49, 5).Value = '  +1.27%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '6.80'
$ws.Range('E50').NumberFormat = '@'
$ws.Cells.Item(50, 5).Value = '  -1.40%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Cells.Item(51, 5).Value = '  -0.23%  '
